# Auto-generated edit script applying the scheduled price-refresh diff
# to the Seraph_Profits workbook (per-sheet profit/cost recalculation columns H:N).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 104.666664  # H33 (was 104.833336)
$ws.Cells.Item(33, 9).Value = 97.59999999999999  # I33 (was 97.8)
$ws.Cells.Item(33, 11).Value = 97.59999999999999  # K33 (was 97.8)
$ws.Cells.Item(33, 13).Value = 131.4  # M33 (was 131.2)
$ws.Cells.Item(63, 8).Value = 50000  # H63 (was 0)
$ws.Cells.Item(63, 10).Value = 50000  # J63 (was 0)
$ws.Cells.Item(63, 12).Value = 50000  # L63 (was 0)
$ws.Cells.Item(63, 14).Value = -51248  # N63 (was None)
$ws.Cells.Item(64, 8).Value = 3519.6  # H64 (was 3999.25)
$ws.Cells.Item(64, 9).Value = 3749.5  # I64 (was 3999.25)
$ws.Cells.Item(64, 10).Value = 2600  # J64 (was 0)
$ws.Cells.Item(64, 11).Value = 3749.5  # K64 (was 3999.25)
$ws.Cells.Item(64, 12).Value = 2600  # L64 (was 0)
$ws.Cells.Item(64, 13).Value = -3501.5  # M64 (was -3751.25)
$ws.Cells.Item(64, 14).Value = -3096  # N64 (was None)
$ws.Cells.Item(66, 8).Value = 50000  # H66 (was 0)
$ws.Cells.Item(66, 10).Value = 50000  # J66 (was 0)
$ws.Cells.Item(66, 12).Value = 150000  # L66 (was 0)
$ws.Cells.Item(66, 14).Value = -156240  # N66 (was None)
$ws.Cells.Item(67, 8).Value = 3519.6  # H67 (was 3999.25)
$ws.Cells.Item(67, 9).Value = 3749.5  # I67 (was 3999.25)
$ws.Cells.Item(67, 10).Value = 2600  # J67 (was 0)
$ws.Cells.Item(67, 11).Value = 3749.5  # K67 (was 3999.25)
$ws.Cells.Item(67, 12).Value = 2600  # L67 (was 0)
$ws.Cells.Item(67, 13).Value = -2891.5  # M67 (was -3141.25)
$ws.Cells.Item(67, 14).Value = -4316  # N67 (was None)
$ws.Cells.Item(70, 8).Value = 5946.1924  # H70 (was 6527.3184)
$ws.Cells.Item(70, 9).Value = 6682.2354  # I70 (was 6974.875)
$ws.Cells.Item(70, 10).Value = 4555.8887  # J70 (was 5333.8335)
$ws.Cells.Item(70, 11).Value = 20046.7062  # K70 (was 20924.625)
$ws.Cells.Item(70, 12).Value = 13667.6661  # L70 (was 16001.5005)
$ws.Cells.Item(70, 13).Value = -19776.7062  # M70 (was -20654.625)
$ws.Cells.Item(70, 14).Value = -14207.6661  # N70 (was -16541.5005)
$ws.Cells.Item(73, 8).Value = 5946.1924  # H73 (was 6527.3184)
$ws.Cells.Item(73, 9).Value = 6682.2354  # I73 (was 6974.875)
$ws.Cells.Item(73, 10).Value = 4555.8887  # J73 (was 5333.8335)
$ws.Cells.Item(73, 11).Value = 20046.7062  # K73 (was 20924.625)
$ws.Cells.Item(73, 12).Value = 13667.6661  # L73 (was 16001.5005)
$ws.Cells.Item(73, 13).Value = -19110.7062  # M73 (was -19988.625)
$ws.Cells.Item(73, 14).Value = -15539.6661  # N73 (was -17873.5005)
$ws.Cells.Item(98, 8).Value = 1156  # H98 (was 679.9474)
$ws.Cells.Item(98, 9).Value = 1156  # I98 (was 679.9474)
$ws.Cells.Item(98, 11).Value = 1156  # K98 (was 679.9474)
$ws.Cells.Item(98, 13).Value = 342  # M98 (was 818.0526)
$ws.Cells.Item(100, 8).Value = 0  # H100 (was 304)
$ws.Cells.Item(100, 9).Value = 0  # I100 (was 304)
$ws.Cells.Item(100, 11).Value = 0  # K100 (was 304)
$ws.Cells.Item(100, 13).ClearContents()  # M100 (was 237)
$ws.Cells.Item(107, 8).Value = 2472  # H107 (was 2514.7144)
$ws.Cells.Item(107, 9).Value = 2472  # I107 (was 2514.7144)
$ws.Cells.Item(107, 11).Value = 2472  # K107 (was 2514.7144)
$ws.Cells.Item(107, 13).Value = -552  # M107 (was -594.7143999999998)
$ws.Cells.Item(111, 8).Value = 1775  # H111 (was 2949.5)
$ws.Cells.Item(111, 9).Value = 700  # I111 (was 850)
$ws.Cells.Item(111, 10).Value = 5000  # J111 (was 3999.25)
$ws.Cells.Item(111, 11).Value = 2100  # K111 (was 2550)
$ws.Cells.Item(111, 12).Value = 15000  # L111 (was 11997.75)
$ws.Cells.Item(111, 13).Value = 967  # M111 (was 517)
$ws.Cells.Item(111, 14).Value = -21134  # N111 (was -18131.75)
$ws.Cells.Item(116, 8).Value = 5450  # H116 (was 7000)
$ws.Cells.Item(116, 9).Value = 2350  # I116 (was 0)
$ws.Cells.Item(116, 11).Value = 2350  # K116 (was 0)
$ws.Cells.Item(116, 13).Value = 1092  # M116 (was None)
$ws.Cells.Item(122, 8).Value = 1156  # H122 (was 679.9474)
$ws.Cells.Item(122, 9).Value = 1156  # I122 (was 679.9474)
$ws.Cells.Item(122, 11).Value = 3468  # K122 (was 2039.8422)
$ws.Cells.Item(122, 13).Value = -1018  # M122 (was 410.1578)
$ws.Cells.Item(132, 8).Value = 2358.1428  # H132 (was 2567.7368)
$ws.Cells.Item(132, 9).Value = 1251.2222  # I132 (was 1361.75)
$ws.Cells.Item(132, 11).Value = 3753.6666  # K132 (was 4085.25)
$ws.Cells.Item(132, 13).Value = -1223.6666  # M132 (was -1555.25)
$ws.Cells.Item(138, 8).Value = 3465.2593  # H138 (was 3498.7778)
$ws.Cells.Item(138, 9).Value = 3982.1667  # I138 (was 4235.636)
$ws.Cells.Item(138, 10).Value = 3051.7334  # J138 (was 2992.1875)
$ws.Cells.Item(138, 11).Value = 11946.5001  # K138 (was 12706.908)
$ws.Cells.Item(138, 12).Value = 9155.200199999999  # L138 (was 8976.5625)
$ws.Cells.Item(138, 13).Value = -6806.500100000001  # M138 (was -7566.908000000001)
$ws.Cells.Item(138, 14).Value = -19435.2002  # N138 (was -19256.5625)

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 4968.5  # H63 (was 4918.857)
$ws.Cells.Item(63, 9).Value = 1194.875  # I63 (was 1108)
$ws.Cells.Item(63, 11).Value = 1194.875  # K63 (was 1108)
$ws.Cells.Item(63, 13).Value = -508.875  # M63 (was -422)
$ws.Cells.Item(66, 8).Value = 4968.5  # H66 (was 4918.857)
$ws.Cells.Item(66, 9).Value = 1194.875  # I66 (was 1108)
$ws.Cells.Item(66, 11).Value = 5974.375  # K66 (was 5540)
$ws.Cells.Item(66, 13).Value = -2542.375  # M66 (was -2108)
$ws.Cells.Item(103, 8).Value = 21032.666  # H103 (was 22375)
$ws.Cells.Item(103, 10).Value = 21032.666  # J103 (was 22375)
$ws.Cells.Item(103, 12).Value = 21032.666  # L103 (was 22375)
$ws.Cells.Item(103, 14).Value = -23376.666  # N103 (was -24719)

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2041.0476  # H105 (was 1993.2273)
$ws.Cells.Item(105, 9).Value = 1335.375  # I105 (was 1315)
$ws.Cells.Item(105, 11).Value = 1335.375  # K105 (was 1315)
$ws.Cells.Item(105, 13).Value = 411.625  # M105 (was 432)
$ws.Cells.Item(134, 8).Value = 2896.2222  # H134 (was 2984.5)
$ws.Cells.Item(134, 9).Value = 2727  # I134 (was 2803.7144)
$ws.Cells.Item(134, 11).Value = 8181  # K134 (was 8411.143199999999)
$ws.Cells.Item(134, 13).Value = -5646  # M134 (was -5876.143199999999)

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 15051  # H4 (was 30002)
$ws.Cells.Item(4, 10).Value = 15051  # J4 (was 30002)
$ws.Cells.Item(4, 12).Value = 15051  # L4 (was 30002)
$ws.Cells.Item(4, 14).Value = -15275  # N4 (was -30226)
$ws.Cells.Item(31, 8).Value = 5246.25  # H31 (was 4905.091)
$ws.Cells.Item(31, 9).Value = 3772.2856  # I31 (was 3425.75)
$ws.Cells.Item(31, 10).Value = 7309.8  # J31 (was 8850)
$ws.Cells.Item(31, 11).Value = 3772.2856  # K31 (was 3425.75)
$ws.Cells.Item(31, 12).Value = 7309.8  # L31 (was 8850)
$ws.Cells.Item(31, 13).Value = -3477.2856  # M31 (was -3130.75)
$ws.Cells.Item(31, 14).Value = -7899.8  # N31 (was -9440)
$ws.Cells.Item(34, 8).Value = 5246.25  # H34 (was 4905.091)
$ws.Cells.Item(34, 9).Value = 3772.2856  # I34 (was 3425.75)
$ws.Cells.Item(34, 10).Value = 7309.8  # J34 (was 8850)
$ws.Cells.Item(34, 11).Value = 3772.2856  # K34 (was 3425.75)
$ws.Cells.Item(34, 12).Value = 7309.8  # L34 (was 8850)
$ws.Cells.Item(34, 13).Value = -3570.2856  # M34 (was -3223.75)
$ws.Cells.Item(34, 14).Value = -7713.8  # N34 (was -9254)
$ws.Cells.Item(59, 8).Value = 33941.668  # H59 (was 35255.555)
$ws.Cells.Item(59, 10).Value = 34080  # J59 (was 35828.57)
$ws.Cells.Item(59, 12).Value = 34080  # L59 (was 35828.57)
$ws.Cells.Item(59, 14).Value = -36370  # N59 (was -38118.57)
$ws.Cells.Item(132, 8).Value = 3904.5417  # H132 (was 3832.1072)
$ws.Cells.Item(132, 9).Value = 2127.4  # I132 (was 2219.2222)
$ws.Cells.Item(132, 10).Value = 6866.4443  # J132 (was 6735.3)
$ws.Cells.Item(132, 11).Value = 6382.200000000001  # K132 (was 6657.6666)
$ws.Cells.Item(132, 12).Value = 20599.3329  # L132 (was 20205.9)
$ws.Cells.Item(132, 13).Value = -3852.200000000001  # M132 (was -4127.6666)
$ws.Cells.Item(132, 14).Value = -25659.3329  # N132 (was -25265.9)

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(140, 8).Value = 4183.1  # H140 (was 4124.5)
$ws.Cells.Item(140, 9).Value = 3366.6  # I140 (was 3833.3333)
$ws.Cells.Item(140, 10).Value = 4999.6  # J140 (was 4998)
$ws.Cells.Item(140, 11).Value = 10099.8  # K140 (was 11499.9999)
$ws.Cells.Item(140, 12).Value = 14998.8  # L140 (was 14994)
$ws.Cells.Item(140, 13).Value = -4919.799999999999  # M140 (was -6319.999899999999)
$ws.Cells.Item(140, 14).Value = -25358.8  # N140 (was -25354)

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 2238  # H5 (was 2675)
$ws.Cells.Item(5, 9).Value = 2172.5  # I5 (was 2733.3333)
$ws.Cells.Item(5, 11).Value = 2172.5  # K5 (was 2733.3333)
$ws.Cells.Item(5, 13).Value = -2060.5  # M5 (was -2621.3333)
$ws.Cells.Item(6, 8).Value = 1902.4  # H6 (was 0)
$ws.Cells.Item(6, 10).Value = 1902.4  # J6 (was 0)
$ws.Cells.Item(6, 12).Value = 1902.4  # L6 (was 0)
$ws.Cells.Item(6, 14).Value = -2128.4  # N6 (was None)
$ws.Cells.Item(16, 8).Value = 1902.4  # H16 (was 0)
$ws.Cells.Item(16, 10).Value = 1902.4  # J16 (was 0)
$ws.Cells.Item(16, 12).Value = 1902.4  # L16 (was 0)
$ws.Cells.Item(16, 14).Value = -2402.4  # N16 (was None)
$ws.Cells.Item(39, 8).Value = 50001  # H39 (was 55261)
$ws.Cells.Item(39, 10).Value = 50001  # J39 (was 55261)
$ws.Cells.Item(39, 12).Value = 50001  # L39 (was 55261)
$ws.Cells.Item(39, 14).Value = -51065  # N39 (was -56325)
$ws.Cells.Item(69, 8).Value = 45000  # H69 (was 0)
$ws.Cells.Item(69, 10).Value = 45000  # J69 (was 0)
$ws.Cells.Item(69, 12).Value = 45000  # L69 (was 0)
$ws.Cells.Item(69, 14).Value = -46498  # N69 (was None)
$ws.Cells.Item(72, 8).Value = 45000  # H72 (was 0)
$ws.Cells.Item(72, 10).Value = 45000  # J72 (was 0)
$ws.Cells.Item(72, 12).Value = 135000  # L72 (was 0)
$ws.Cells.Item(72, 14).Value = -142488  # N72 (was None)
$ws.Cells.Item(113, 8).Value = 12520  # H113 (was 12144)
$ws.Cells.Item(113, 10).Value = 19272.455  # J113 (was 18083.084)
$ws.Cells.Item(113, 12).Value = 19272.455  # L113 (was 18083.084)
$ws.Cells.Item(113, 14).Value = -23612.455  # N113 (was -22423.084)
$ws.Cells.Item(118, 8).Value = 42555.555  # H118 (was 43000)
$ws.Cells.Item(118, 10).Value = 42555.555  # J118 (was 43000)
$ws.Cells.Item(118, 12).Value = 42555.555  # L118 (was 43000)
$ws.Cells.Item(118, 14).Value = -45869.555  # N118 (was -46314)
$ws.Cells.Item(126, 8).Value = 4186  # H126 (was 3383.8096)
$ws.Cells.Item(126, 9).Value = 3702.2  # I126 (was 3041.75)
$ws.Cells.Item(126, 10).Value = 4488.375  # J126 (was 3594.3076)
$ws.Cells.Item(126, 11).Value = 11106.6  # K126 (was 9125.25)
$ws.Cells.Item(126, 12).Value = 13465.125  # L126 (was 10782.9228)
$ws.Cells.Item(126, 13).Value = -8636.599999999999  # M126 (was -6655.25)
$ws.Cells.Item(126, 14).Value = -18405.125  # N126 (was -15722.9228)

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(4, 8).Value = 7377223  # H4 (was 5252.25)
$ws.Cells.Item(4, 9).Value = 11060835  # I4 (was 3669.6667)
$ws.Cells.Item(4, 10).Value = 9999.333000000001  # J4 (was 10000)
$ws.Cells.Item(4, 11).Value = 11060835  # K4 (was 3669.6667)
$ws.Cells.Item(4, 12).Value = 9999.333000000001  # L4 (was 10000)
$ws.Cells.Item(4, 13).Value = -11060722  # M4 (was -3556.6667)
$ws.Cells.Item(4, 14).Value = -10225.333  # N4 (was -10226)
$ws.Cells.Item(5, 8).Value = 40003  # H5 (was 40010)
$ws.Cells.Item(5, 9).Value = 40000  # I5 (was 0)
$ws.Cells.Item(5, 10).Value = 40004.5  # J5 (was 40010)
$ws.Cells.Item(5, 11).Value = 40000  # K5 (was 0)
$ws.Cells.Item(5, 12).Value = 40004.5  # L5 (was 40010)
$ws.Cells.Item(5, 13).Value = -39887  # M5 (was None)
$ws.Cells.Item(5, 14).Value = -40230.5  # N5 (was -40236)
$ws.Cells.Item(16, 8).Value = 9248.833000000001  # H16 (was 9549.125)
$ws.Cells.Item(16, 9).Value = 9138.799999999999  # I16 (was 9115.833000000001)
$ws.Cells.Item(16, 10).Value = 9799  # J16 (was 10849)
$ws.Cells.Item(16, 11).Value = 9138.799999999999  # K16 (was 9115.833000000001)
$ws.Cells.Item(16, 12).Value = 9799  # L16 (was 10849)
$ws.Cells.Item(16, 13).Value = -8968.799999999999  # M16 (was -8945.833000000001)
$ws.Cells.Item(16, 14).Value = -10139  # N16 (was -11189)
$ws.Cells.Item(28, 8).Value = 7377223  # H28 (was 5252.25)
$ws.Cells.Item(28, 9).Value = 11060835  # I28 (was 3669.6667)
$ws.Cells.Item(28, 10).Value = 9999.333000000001  # J28 (was 10000)
$ws.Cells.Item(28, 11).Value = 11060835  # K28 (was 3669.6667)
$ws.Cells.Item(28, 12).Value = 9999.333000000001  # L28 (was 10000)
$ws.Cells.Item(28, 13).Value = -11060603  # M28 (was -3437.6667)
$ws.Cells.Item(28, 14).Value = -10463.333  # N28 (was -10464)
$ws.Cells.Item(37, 8).Value = 7377223  # H37 (was 5252.25)
$ws.Cells.Item(37, 9).Value = 11060835  # I37 (was 3669.6667)
$ws.Cells.Item(37, 10).Value = 9999.333000000001  # J37 (was 10000)
$ws.Cells.Item(37, 11).Value = 11060835  # K37 (was 3669.6667)
$ws.Cells.Item(37, 12).Value = 9999.333000000001  # L37 (was 10000)
$ws.Cells.Item(37, 13).Value = -11060728  # M37 (was -3562.6667)
$ws.Cells.Item(37, 14).Value = -10213.333  # N37 (was -10214)
$ws.Cells.Item(46, 8).Value = 1500  # H46 (was 2625)
$ws.Cells.Item(46, 9).Value = 1500  # I46 (was 2625)
$ws.Cells.Item(46, 11).Value = 1500  # K46 (was 2625)
$ws.Cells.Item(46, 13).Value = -1312  # M46 (was -2437)
$ws.Cells.Item(68, 8).Value = 3747.0952  # H68 (was 3823.3684)
$ws.Cells.Item(68, 9).Value = 3466.0557  # I68 (was 3521.5)
$ws.Cells.Item(68, 11).Value = 3466.0557  # K68 (was 3521.5)
$ws.Cells.Item(68, 13).Value = -2717.0557  # M68 (was -2772.5)
$ws.Cells.Item(71, 8).Value = 3747.0952  # H71 (was 3823.3684)
$ws.Cells.Item(71, 9).Value = 3466.0557  # I71 (was 3521.5)
$ws.Cells.Item(71, 11).Value = 17330.2785  # K71 (was 17607.5)
$ws.Cells.Item(71, 13).Value = -13586.2785  # M71 (was -13863.5)
$ws.Cells.Item(74, 8).Value = 62000  # H74 (was 65000)
$ws.Cells.Item(74, 10).Value = 62000  # J74 (was 65000)
$ws.Cells.Item(74, 12).Value = 62000  # L74 (was 65000)
$ws.Cells.Item(74, 14).Value = -63996  # N74 (was -66996)
$ws.Cells.Item(77, 8).Value = 62000  # H77 (was 65000)
$ws.Cells.Item(77, 10).Value = 62000  # J77 (was 65000)
$ws.Cells.Item(77, 12).Value = 186000  # L77 (was 195000)
$ws.Cells.Item(77, 14).Value = -195984  # N77 (was -204984)
$ws.Cells.Item(82, 8).Value = 144570.72  # H82 (was 112598.89)
$ws.Cells.Item(82, 9).Value = 1999.3334  # I82 (was 1841.5714)
$ws.Cells.Item(82, 10).Value = 999999  # J82 (was 500249.5)
$ws.Cells.Item(82, 11).Value = 1999.3334  # K82 (was 1841.5714)
$ws.Cells.Item(82, 12).Value = 999999  # L82 (was 500249.5)
$ws.Cells.Item(82, 13).Value = -1638.3334  # M82 (was -1480.5714)
$ws.Cells.Item(82, 14).Value = -1000721  # N82 (was -500971.5)
$ws.Cells.Item(85, 8).Value = 144570.72  # H85 (was 112598.89)
$ws.Cells.Item(85, 9).Value = 1999.3334  # I85 (was 1841.5714)
$ws.Cells.Item(85, 10).Value = 999999  # J85 (was 500249.5)
$ws.Cells.Item(85, 11).Value = 1999.3334  # K85 (was 1841.5714)
$ws.Cells.Item(85, 12).Value = 999999  # L85 (was 500249.5)
$ws.Cells.Item(85, 13).Value = -751.3334  # M85 (was -593.5714)
$ws.Cells.Item(85, 14).Value = -1002495  # N85 (was -502745.5)
$ws.Cells.Item(122, 8).Value = 2998  # H122 (was 3133.3333)
$ws.Cells.Item(122, 9).Value = 2998  # I122 (was 3133.3333)
$ws.Cells.Item(122, 11).Value = 8994  # K122 (was 9399.999899999999)
$ws.Cells.Item(122, 13).Value = -6544  # M122 (was -6949.999899999999)

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(8, 8).Value = 2033.6666  # H8 (was 1757.8572)
$ws.Cells.Item(8, 9).Value = 350  # I8 (was 226.5)
$ws.Cells.Item(8, 11).Value = 350  # K8 (was 226.5)
$ws.Cells.Item(8, 13).Value = -210  # M8 (was -86.5)
$ws.Cells.Item(62, 8).Value = 7653.846  # H62 (was 5761.1763)
$ws.Cells.Item(62, 9).Value = 6800  # I62 (was 2691.4285)
$ws.Cells.Item(62, 11).Value = 6800  # K62 (was 2691.4285)
$ws.Cells.Item(62, 13).Value = -6176  # M62 (was -2067.4285)
$ws.Cells.Item(65, 8).Value = 7653.846  # H65 (was 5761.1763)
$ws.Cells.Item(65, 9).Value = 6800  # I65 (was 2691.4285)
$ws.Cells.Item(65, 11).Value = 34000  # K65 (was 13457.1425)
$ws.Cells.Item(65, 13).Value = -30880  # M65 (was -10337.1425)

